# Apply the "newest EPS-US files" revision to
# CO2 Abated per Unit Land Area by Impr For Mgmt.xlsx
#
# Summary of the change:
#  - The "Calculations" worksheet (the intermediate low/high/average
#    calculation) is removed entirely.
#  - On the "CApULAbIFM" worksheet, the CO2 Abated (g) figure is now
#    computed directly with the formula =1.5*10^6 instead of pulling
#    the old Calculations!A6 result.
#  - On the "About" worksheet, the source citation is simplified to
#    "consultation with American Forest Foundation" and the supporting
#    rows (year, low/high estimate text, EPA hyperlink/citation, page
#    reference) are deleted.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsCalcs      = $wb.Worksheets.Item("Calculations")
$wsCApULAbIFM = $wb.Worksheets.Item("CApULAbIFM")
$wsAbout      = $wb.Worksheets.Item("About")

# --- CApULAbIFM sheet: replace the link to Calculations!A6 with a direct formula ---
$wsCApULAbIFM.Range("B2").Formula = "=1.5*10^6"

# --- Remove the now unnecessary Calculations worksheet ---
$wsCalcs.Delete()

# Sheet collection indices shifted after the delete above; re-fetch the
# remaining worksheet references so later calls (Activate, etc.) operate
# on live objects instead of a stale pre-delete handle.
$wsCApULAbIFM = $wb.Worksheets.Item("CApULAbIFM")
$wsAbout      = $wb.Worksheets.Item("About")

# --- About sheet: update the source citation and drop the old supporting rows ---
$wsAbout.Hyperlinks.Delete()
$wsAbout.Range("B3").Value = "consultation with American Forest Foundation"
$wsAbout.Rows("4:7").Delete()
$wsAbout.Range("A5").Font.Bold = $false

# The workbook no longer uses any hyperlink-styled cells, drop the style
$wb.Styles.Item("Hyperlink").Delete()

# Restore the cursor/selection positions left behind in the saved file
$wsCApULAbIFM.Activate()
$wsCApULAbIFM.Range("B3").Select()

$wsAbout.Activate()
$wsAbout.Range("C17").Select()
